# Auto-generated Excel COM-interop script applying Sheets data refresh
# (scheduled runner update to currentAveragePrice / Leve price / profit columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1265.6666
$ws.Range("I19").Value = 1408.1
$ws.Range("J19").Value = 1163.9286
$ws.Range("K19").Value = 1408.1
$ws.Range("L19").Value = 1163.9286
$ws.Range("M19").Value = -1233.1
$ws.Range("N19").Value = -1513.9286

$ws.Range("H70").Value = 1272
$ws.Range("I70").Value = 1196
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 3588
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -3318
$ws.Range("N70").Value = -5040

$ws.Range("H73").Value = 1272
$ws.Range("I73").Value = 1196
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 3588
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -2652
$ws.Range("N73").Value = -6372

$ws.Range("H94").Value = 90910120
$ws.Range("I94").Value = 1130.4
$ws.Range("K94").Value = 1130.4
$ws.Range("M94").Value = -679.4000000000001

$ws.Range("H103").Value = 779.2759
$ws.Range("I103").Value = 600
$ws.Range("J103").Value = 785.6786
$ws.Range("K103").Value = 1800
$ws.Range("L103").Value = 2357.0358
$ws.Range("M103").Value = -1214
$ws.Range("N103").Value = -3529.0358

$ws.Range("H115").Value = 25087.5
$ws.Range("I115").Value = 25087.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 75262.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -73695.5
$ws.Range("N115").ClearContents()

$ws.Range("H127").Value = 2184.5574
$ws.Range("J127").Value = 2336.4465
$ws.Range("L127").Value = 7009.3395
$ws.Range("N127").Value = -16929.3395

$ws.Range("H129").Value = 930.25
$ws.Range("I129").Value = 1000
$ws.Range("J129").Value = 928.7659
$ws.Range("K129").Value = 3000
$ws.Range("L129").Value = 2786.2977
$ws.Range("M129").Value = 2000
$ws.Range("N129").Value = -12786.2977

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 38126.43
$ws.Range("I45").Value = 54582.316
$ws.Range("J45").Value = 3386.2222
$ws.Range("K45").Value = 54582.316
$ws.Range("L45").Value = 3386.2222
$ws.Range("M45").Value = -54205.316
$ws.Range("N45").Value = -4140.2222

$ws.Range("H74").Value = 1347.4857
$ws.Range("I74").Value = 711.5833
$ws.Range("J74").Value = 2734.9092
$ws.Range("K74").Value = 711.5833
$ws.Range("L74").Value = 2734.9092
$ws.Range("M74").Value = 162.4167
$ws.Range("N74").Value = -4482.9092

$ws.Range("H77").Value = 1347.4857
$ws.Range("I77").Value = 711.5833
$ws.Range("J77").Value = 2734.9092
$ws.Range("K77").Value = 3557.9165
$ws.Range("L77").Value = 13674.546
$ws.Range("M77").Value = 810.0834999999997
$ws.Range("N77").Value = -22410.546

$ws.Range("H102").Value = 58596.832
$ws.Range("I102").Value = 85948.586
$ws.Range("J102").Value = 3893.3333
$ws.Range("K102").Value = 85948.586
$ws.Range("L102").Value = 3893.3333
$ws.Range("M102").Value = -84326.586
$ws.Range("N102").Value = -7137.3333

$ws.Range("H122").Value = 2650.5
$ws.Range("I122").Value = 2743.4285
$ws.Range("K122").Value = 8230.2855
$ws.Range("M122").Value = -5780.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H99").Value = 1141.4706
$ws.Range("I99").Value = 887.5
$ws.Range("J99").Value = 2326.6667
$ws.Range("K99").Value = 887.5
$ws.Range("L99").Value = 2326.6667
$ws.Range("M99").Value = 610.5
$ws.Range("N99").Value = -5322.6667

$ws.Range("H105").Value = 119658.88
$ws.Range("I105").Value = 168929.83
$ws.Range("J105").Value = 92783.82000000001
$ws.Range("K105").Value = 168929.83
$ws.Range("L105").Value = 92783.82000000001
$ws.Range("M105").Value = -167182.83
$ws.Range("N105").Value = -96277.82000000001

$ws.Range("H107").Value = 76959200
$ws.Range("I107").Value = 125057810
$ws.Range("J107").Value = 1439.8
$ws.Range("K107").Value = 125057810
$ws.Range("L107").Value = 1439.8
$ws.Range("M107").Value = -125055890
$ws.Range("N107").Value = -5279.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6946570
$ws.Range("I62").Value = 18519688
$ws.Range("J62").Value = 2699.4
$ws.Range("K62").Value = 18519688
$ws.Range("L62").Value = 2699.4
$ws.Range("M62").Value = -18519064
$ws.Range("N62").Value = -3947.4

$ws.Range("H65").Value = 6946570
$ws.Range("I65").Value = 18519688
$ws.Range("J65").Value = 2699.4
$ws.Range("K65").Value = 92598440
$ws.Range("L65").Value = 13497
$ws.Range("M65").Value = -92595320
$ws.Range("N65").Value = -19737

$ws.Range("H86").Value = 1877.6296
$ws.Range("I86").Value = 1434.7059
$ws.Range("J86").Value = 2630.6
$ws.Range("K86").Value = 1434.7059
$ws.Range("L86").Value = 2630.6
$ws.Range("M86").Value = -311.7058999999999
$ws.Range("N86").Value = -4876.6

$ws.Range("H89").Value = 1877.6296
$ws.Range("I89").Value = 1434.7059
$ws.Range("J89").Value = 2630.6
$ws.Range("K89").Value = 7173.5295
$ws.Range("L89").Value = 13153
$ws.Range("M89").Value = -1557.5295
$ws.Range("N89").Value = -24385

$ws.Range("H96").Value = 22155.75
$ws.Range("J96").Value = 22155.75
$ws.Range("L96").Value = 22155.75
$ws.Range("N96").Value = -27647.75

$ws.Range("H107").Value = 1216.5238
$ws.Range("I107").Value = 1328.4375
$ws.Range("J107").Value = 858.4
$ws.Range("K107").Value = 1328.4375
$ws.Range("L107").Value = 858.4
$ws.Range("M107").Value = 591.5625
$ws.Range("N107").Value = -4698.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2008.2963
$ws.Range("I75").Value = 528
$ws.Range("J75").Value = 2631.5789
$ws.Range("K75").Value = 1584
$ws.Range("L75").Value = 7894.736699999999
$ws.Range("M75").Value = -586
$ws.Range("N75").Value = -9890.736699999999

$ws.Range("H78").Value = 2008.2963
$ws.Range("I78").Value = 528
$ws.Range("J78").Value = 2631.5789
$ws.Range("K78").Value = 4752
$ws.Range("L78").Value = 23684.2101
$ws.Range("M78").Value = 240
$ws.Range("N78").Value = -33668.2101

$ws.Range("H97").Value = 666.8333
$ws.Range("J97").Value = 699.6
$ws.Range("L97").Value = 2098.8
$ws.Range("N97").Value = -3090.8

$ws.Range("H98").Value = 66551
$ws.Range("I98").Value = 1121.2
$ws.Range("J98").Value = 93813.414
$ws.Range("K98").Value = 3363.6
$ws.Range("L98").Value = 281440.242
$ws.Range("M98").Value = -1865.6
$ws.Range("N98").Value = -284436.242

$ws.Range("H131").Value = 858.1414
$ws.Range("J131").Value = 858.1414
$ws.Range("L131").Value = 2574.4242
$ws.Range("N131").Value = -12654.4242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1876
$ws.Range("I113").Value = 1920.3334
$ws.Range("J113").Value = 1859.375
$ws.Range("K113").Value = 1920.3334
$ws.Range("L113").Value = 1859.375
$ws.Range("M113").Value = 249.6666
$ws.Range("N113").Value = -6199.375

$ws.Range("H120").Value = 34298.25
$ws.Range("J120").Value = 34298.25
$ws.Range("L120").Value = 34298.25
$ws.Range("N120").Value = -43974.25

$ws.Range("H122").Value = 1166.5
$ws.Range("I122").Value = 1180
$ws.Range("J122").Value = 1099
$ws.Range("K122").Value = 3540
$ws.Range("L122").Value = 3297
$ws.Range("M122").Value = -1090
$ws.Range("N122").Value = -8197

$ws.Range("H126").Value = 3606.182
$ws.Range("I126").Value = 3203.5
$ws.Range("J126").Value = 4089.4
$ws.Range("K126").Value = 9610.5
$ws.Range("L126").Value = 12268.2
$ws.Range("M126").Value = -7140.5
$ws.Range("N126").Value = -17208.2

$ws.Range("H134").Value = 24616.125
$ws.Range("J134").Value = 24616.125
$ws.Range("L134").Value = 73848.375
$ws.Range("N134").Value = -78918.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2878.9333
$ws.Range("I61").Value = 2320.4
$ws.Range("J61").Value = 3996
$ws.Range("K61").Value = 2320.4
$ws.Range("L61").Value = 3996
$ws.Range("M61").Value = -2118.4
$ws.Range("N61").Value = -4400

$ws.Range("H100").Value = 1372.3334
$ws.Range("I100").Value = 1384
$ws.Range("J100").Value = 1349
$ws.Range("K100").Value = 1384
$ws.Range("L100").Value = 1349
$ws.Range("M100").Value = -843
$ws.Range("N100").Value = -2431

$ws.Range("H113").Value = 2878.9333
$ws.Range("I113").Value = 2320.4
$ws.Range("J113").Value = 3996
$ws.Range("K113").Value = 2320.4
$ws.Range("L113").Value = 3996
$ws.Range("M113").Value = -150.4000000000001
$ws.Range("N113").Value = -8336

$ws.Range("H136").Value = 1356.3143
$ws.Range("I136").Value = 1302.3572
$ws.Range("J136").Value = 1572.1428
$ws.Range("K136").Value = 3907.0716
$ws.Range("L136").Value = 4716.428400000001
$ws.Range("M136").Value = -1357.0716
$ws.Range("N136").Value = -9816.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 142858600
$ws.Range("I96").Value = 166668160
$ws.Range("J96").Value = 1200
$ws.Range("K96").Value = 166668160
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = -166666787
$ws.Range("N96").Value = -3946

$ws.Range("H121").Value = 39999
$ws.Range("J121").Value = 39999
$ws.Range("L121").Value = 39999
$ws.Range("N121").Value = -43493

$ws.Range("H122").Value = 1900.4
$ws.Range("I122").Value = 999.8333
$ws.Range("J122").Value = 3251.25
$ws.Range("K122").Value = 2999.4999
$ws.Range("L122").Value = 9753.75
$ws.Range("M122").Value = -549.4998999999998
$ws.Range("N122").Value = -14653.75

$ws.Range("H126").Value = 1452.6923
$ws.Range("I126").Value = 1452.6923
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4358.0769
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1888.0769
$ws.Range("N126").ClearContents()
